# "Generate Report for Handback"
# The handback for the two localized files (d44ffe7d... and e1a0821e...) has
# completed. Update the localization-status report so that:
#   - the per-language "Status" cells on the Overview sheet read as handed
#     back (and in sync with en-US) instead of "Ready for handoff"
#   - each language sheet (zh-cn, de-de) gets the "Latest Target File" /
#     "Latest Handback File" filled in (with a hyperlink on the target file,
#     matching the existing source-file hyperlink) and a fresh
#     "Latest Handback DateTime" stamp

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status cells for both rows.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.15
$overview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------
# Helper data: per language-sheet, per row, the handback info.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Capture the existing hyperlink addresses for column A (Source File Name)
# so the new column I (Latest Target File) hyperlinks point at the same
# GitHub blob URLs.
$zhcnAddrs = @()
foreach ($hl in $zhcn.Hyperlinks) { $zhcnAddrs += $hl.Address }

$dedeAddrs = @()
foreach ($hl in $dede.Hyperlinks) { $dedeAddrs += $hl.Address }

$file1Name = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.md"
$file2Name = "e1a0821e-2950-410f-ac19-156cb9e5b724.md"

# ---------------------------------------------------------------------
# zh-cn sheet: rows 2 (d44ffe7d...) and 3 (e1a0821e...)
# ---------------------------------------------------------------------
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

$zhcn.Range("J2").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-18 14:51:55"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $zhcnAddrs[0], "", "", $file1Name)

$zhcn.Range("J3").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-18 14:51:55"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $zhcnAddrs[1], "", "", $file2Name)

$zhcn.Columns.Item(3).ColumnWidth = 29.15
$zhcn.Columns.Item(9).ColumnWidth = 39.15
$zhcn.Columns.Item(10).ColumnWidth = 39.15

# ---------------------------------------------------------------------
# de-de sheet: rows 2 (d44ffe7d...) and 3 (e1a0821e...)
# ---------------------------------------------------------------------
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

$dede.Range("J2").Value = "d44ffe7d-90fc-4235-9238-eb4b6785fa30.0341931d9c5303e347d81090f2db43f531e92132.de-de.xlf"
$dede.Range("K2").Value = "2016-08-18 14:52:12"
$dede.Hyperlinks.Add($dede.Range("I2"), $dedeAddrs[0], "", "", $file1Name)

$dede.Range("J3").Value = "e1a0821e-2950-410f-ac19-156cb9e5b724.17a387e0c195beaa3ef55d3ec3d993803b4a9f72.de-de.xlf"
$dede.Range("K3").Value = "2016-08-18 14:52:12"
$dede.Hyperlinks.Add($dede.Range("I3"), $dedeAddrs[1], "", "", $file2Name)

$dede.Columns.Item(3).ColumnWidth = 29.15
$dede.Columns.Item(9).ColumnWidth = 39.15
$dede.Columns.Item(10).ColumnWidth = 39.15

Write-Host "Handback report generated."
